$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new test rows were inserted above the existing "andrian" / "andrian2"
# rows (transactional insert / notification-import test data), pushing the
# previous row 2-3 data down to row 4-5.
$ws.Rows("2:3").Insert()

# New row 2 - "aldo" test user
$ws.Range("A2").Value = "aldo"
$ws.Range("B2").Value = "asd"
$ws.Range("C2").Value = "asdasd"
$ws.Range("D2").Value = "asdf"
$ws.Range("E2").Value = "anasdkajsdh"
$ws.Range("F2").Value = "asdf@sad"
$ws.Range("G2").Value = "admin"

# New row 3 - "andrianasd2" test user
$ws.Range("A3").Value = "andrianasd2"
$ws.Range("B3").Value = "andriasdn2"
$ws.Range("C3").Value = "dasdf"
$ws.Range("D3").Value = 9871
$ws.Range("E3").Value = "dasd"
$ws.Range("F3").Value = "adaff@ad"
$ws.Range("G3").Value = "lecturer"

# Match the saved selection state from the edit
$ws.Range("F3").Select()
